$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'67.731.58"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +1.25%  '
$ws.Range("E2").ClearFormats()

# Row 3
$ws.Range("D3").Value = "'2.620.20"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +0.86%  '
$ws.Range("E3").ClearFormats()

# Row 4
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("E4").ClearFormats()

# Row 5
$ws.Range("D5").Value = "'604.42"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +1.89%  '
$ws.Range("E5").ClearFormats()

# Row 6
$ws.Range("D6").Value = "'155.15"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.88%  '
$ws.Range("E6").ClearFormats()

# Row 7
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("E7").ClearFormats()

# Row 8
$ws.Range("E8").Value = '  +1.52%  '
$ws.Range("E8").ClearFormats()

# Row 9
$ws.Range("D9").Value = "'2.621.64"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +0.92%  '
$ws.Range("E9").ClearFormats()

# Row 10
$ws.Range("D10").Value = "'0.125"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +7.26%  '
$ws.Range("E10").ClearFormats()

# Row 11
$ws.Range("E11").Value = '  +0.97%  '
$ws.Range("E11").ClearFormats()

# Row 12
$ws.Range("D12").Value = "'5.25"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +0.61%  '
$ws.Range("E12").ClearFormats()

# Row 13
$ws.Range("E13").Value = '  -1.51%  '
$ws.Range("E13").ClearFormats()

# Row 14
$ws.Range("D14").Value = "'28.08"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +0.62%  '
$ws.Range("E14").ClearFormats()

# Row 15
$ws.Range("D15").Value = "'0.0000185"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +2.16%  '
$ws.Range("E15").ClearFormats()

# Row 16
$ws.Range("D16").Value = "'3.102.01"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +1.11%  '
$ws.Range("E16").ClearFormats()

# Row 17
$ws.Range("D17").Value = "'67.606.99"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +1.40%  '
$ws.Range("E17").ClearFormats()

# Row 18
$ws.Range("D18").Value = "'2.620.93"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +0.94%  '
$ws.Range("E18").ClearFormats()

# Row 19
$ws.Range("B19").Value = 'BitcoinCash'
$ws.Range("B19").ClearFormats()
$ws.Range("C19").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("C19").ClearFormats()
$ws.Range("D19").Value = "'371.03"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +4.02%  '
$ws.Range("E19").ClearFormats()

# Row 20
$ws.Range("B20").Value = 'Chainlink'
$ws.Range("B20").ClearFormats()
$ws.Range("C20").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("C20").ClearFormats()
$ws.Range("D20").Value = "'11.26"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -0.54%  '
$ws.Range("E20").ClearFormats()

# Row 21
$ws.Range("D21").Value = "'7.59"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -2.21%  '
$ws.Range("E21").ClearFormats()

# Row 22
$ws.Range("E22").Value = '  -0.48%  '
$ws.Range("E22").ClearFormats()

# Row 23
$ws.Range("D23").Value = "'2.12"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +5.02%  '
$ws.Range("E23").ClearFormats()

# Row 24
$ws.Range("E24").Value = '  +0.06%  '
$ws.Range("E24").ClearFormats()

# Row 25
$ws.Range("D25").Value = "'70.38"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +6.28%  '
$ws.Range("E25").ClearFormats()

# Row 26
$ws.Range("E26").Value = '  -2.40%  '
$ws.Range("E26").ClearFormats()

# Row 27
$ws.Range("D27").Value = "'0.0000105"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +1.48%  '
$ws.Range("E27").ClearFormats()

# Row 28
$ws.Range("D28").Value = "'2.751.99"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +1.03%  '
$ws.Range("E28").ClearFormats()

# Row 29
$ws.Range("B29").Value = 'Binance-PegBSC-USD'
$ws.Range("B29").ClearFormats()
$ws.Range("C29").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("C29").ClearFormats()
$ws.Range("D29").Value = "'1.01"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +0.71%  '
$ws.Range("E29").ClearFormats()

# Row 30
$ws.Range("B30").Value = 'Bittensor'
$ws.Range("B30").ClearFormats()
$ws.Range("C30").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("C30").ClearFormats()
$ws.Range("D30").Value = "'584.51"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -1.90%  '
$ws.Range("E30").ClearFormats()

# Row 31
$ws.Range("D31").Value = "'1.44"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -0.39%  '
$ws.Range("E31").ClearFormats()

# Row 32
$ws.Range("D32").Value = "'7.94"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -0.94%  '
$ws.Range("E32").ClearFormats()

# Row 33
$ws.Range("E33").Value = '  -0.12%  '
$ws.Range("E33").ClearFormats()

# Row 34
$ws.Range("E34").Value = '  -2.37%  '
$ws.Range("E34").ClearFormats()

# Row 35
$ws.Range("E35").Value = '  -0.01%  '
$ws.Range("E35").ClearFormats()

# Row 36
$ws.Range("E36").Value = '  -1.65%  '
$ws.Range("E36").ClearFormats()

# Row 37
$ws.Range("D37").Value = "'4.99"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +0.30%  '
$ws.Range("E37").ClearFormats()

# Row 38
$ws.Range("D38").Value = "'19.50"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +0.33%  '
$ws.Range("E38").ClearFormats()

# Row 39
$ws.Range("D39").Value = "'157.11"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +1.81%  '
$ws.Range("E39").ClearFormats()

# Row 40
$ws.Range("D40").Value = "'0.372"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +0.69%  '
$ws.Range("E40").ClearFormats()

# Row 41
$ws.Range("E41").Value = '  -0.83%  '
$ws.Range("E41").ClearFormats()

# Row 42
$ws.Range("D42").Value = "'1.86"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +3.21%  '
$ws.Range("E42").ClearFormats()

# Row 43
$ws.Range("E43").Value = '  +1.06%  '
$ws.Range("E43").ClearFormats()

# Row 44
$ws.Range("D44").Value = "'41.13"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -0.54%  '
$ws.Range("E44").ClearFormats()

# Row 45
$ws.Range("B45").Value = 'WhiteBITCoin'
$ws.Range("B45").ClearFormats()
$ws.Range("C45").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("C45").ClearFormats()
$ws.Range("D45").Value = "'16.43"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -0.16%  '
$ws.Range("E45").ClearFormats()

# Row 46
$ws.Range("B46").Value = 'USDe'
$ws.Range("B46").ClearFormats()
$ws.Range("C46").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("C46").ClearFormats()
$ws.Range("D46").Value = "'0.999"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +0.02%  '
$ws.Range("E46").ClearFormats()

# Row 47
$ws.Range("D47").Value = "'156.29"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +1.21%  '
$ws.Range("E47").ClearFormats()

# Row 48
$ws.Range("D48").Value = "'0.0₆0288"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -5.78%  '
$ws.Range("E48").ClearFormats()

# Row 49
$ws.Range("E49").Value = '  -0.42%  '
$ws.Range("E49").ClearFormats()

# Row 50
$ws.Range("D50").Value = "'21.06"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -0.21%  '
$ws.Range("E50").ClearFormats()

# Row 51
$ws.Range("D51").Value = "'0.627"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +0.69%  '
$ws.Range("E51").ClearFormats()
